$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 'Desodin'
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 'Desodin 60ml Syrup'
$ws.Range("E2").Value = '60 ml'

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 'Dinafex'
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 'Dinafex 60mg Tablet'
$ws.Range("E3").Value = '30''s'

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 'Dinafex'
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 'Dinafex 180mg Tablet'
$ws.Range("E4").Value = '30''s'

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 'Dinafex'
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 'Dinafex 120mg Tablet'
$ws.Range("E5").Value = '30''s'

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 'Dorenta'
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 'Dorenta 50mg Tablet'
$ws.Range("E6").Value = '50''s'

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 'Etorix'
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 'Etorix 120mg Tablet'
$ws.Range("E7").Value = '20''s'

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Etorix'
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 'Etorix 90mg Tablet'
$ws.Range("E8").Value = '30''s'

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'Etorix'
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 'Etorix 60mg Tablet - 40''s'
$ws.Range("E9").Value = '40''s'

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'Fenobac'
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 'Fenobac 100ml Syrup'
$ws.Range("E10").Value = '100ml'

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'Flucloxin'
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 'Flucloxin 500mg Capsule - 36''s'
$ws.Range("E11").Value = '36 ''s'

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 'Flucloxin'
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 'Flucloxin 500mg Capsule'
$ws.Range("E12").Value = '30 ''s'

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 'Geminox'
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 'Geminox 320mg Tablet - 8''s'
$ws.Range("E13").Value = '8 ''s'

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 'Ketonic'
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 'Ketonic 10mg Tablet'
$ws.Range("E14").Value = '20''s'

$ws.Range("A15").Value = 11
$ws.Range("B15").Value = 'Ketonic'
$ws.Range("C15").Value = 14
$ws.Range("D15").Value = 'Ketonic 30mg Injection'
$ws.Range("E15").Value = '5 ''s'

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 'Ketonic'
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 'Ketonic 30mg IM/IV Injection - 4''s'
$ws.Range("E16").Value = '4''s'

$ws.Range("A17").Value = 12
$ws.Range("B17").Value = 'Kynol'
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 'Kynol D 25mg Tablet'
$ws.Range("E17").Value = '60 ''s'

$ws.Range("A18").Value = 12
$ws.Range("B18").Value = 'Kynol'
$ws.Range("C18").Value = 17
$ws.Range("D18").Value = 'Kynol TR 200mg Capsule'
$ws.Range("E18").Value = '30 ''s'

$ws.Range("A19").Value = 12
$ws.Range("B19").Value = 'Kynol'
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 'Kynol TR 100mg Capsule'
$ws.Range("E19").Value = '50 ''s'

$ws.Range("A20").Value = 17
$ws.Range("B20").Value = 'Naprox'
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 'Naprox Plus 500mg Tablet - 30''s'
$ws.Range("E20").Value = '30 ''s'

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 'Oradin'
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 'Oradin Plus Tablet - 40''s'
$ws.Range("E21").Value = '40 ''s'

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 'Osticare'
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 'Osticare Tablet 24''s'
$ws.Range("E22").Value = '24''s'

$ws.Range("A23").Value = 23
$ws.Range("B23").Value = 'Rupaday'
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 'Rupaday Oral Solution 60ml'
$ws.Range("E23").Value = '1''s'

$ws.Range("A24").Value = 24
$ws.Range("B24").Value = 'Sk-Mox'
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 'Sk-Mox 500mg Capsule'
$ws.Range("E24").Value = '48 ''s'

$ws.Range("A25").Value = 35
$ws.Range("B25").Value = 'Zithrox'
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 'Zithrox 15ml Suspension'
$ws.Range("E25").Value = '15 ml'

$ws.Range("A26").Value = 35
$ws.Range("B26").Value = 'Zithrox'
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 'Zithrox 30ml Dry Suspension'
$ws.Range("E26").Value = '30ml'

$ws.Range("A27").Value = 35
$ws.Range("B27").Value = 'Zithrox'
$ws.Range("C27").Value = 26
$ws.Range("D27").Value = 'Zithrox 500mg Tablet'
$ws.Range("E27").Value = '6 ''s'

$ws.Range("A28").Value = 35
$ws.Range("B28").Value = 'Zithrox'
$ws.Range("C28").Value = 27
$ws.Range("D28").Value = 'Zithrox 250mg Tablet - 6''s'
$ws.Range("E28").Value = '6''s'

$ws.Rows.Item(29).Delete()

